$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.582.36"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "2.367.79"
$ws.Range("E3").Value = "  -4.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.12"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -2.55%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.35"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -6.53%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -4.16%  "
$ws.Range("E8").Value = "  +0.04%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -4.24%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -2.31%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.21"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -8.69%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.109"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "2.732.53"
$ws.Range("E13").Value = "  -4.17%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.53"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -5.21%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.03"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "2.377.08"
$ws.Range("E16").Value = "  -3.52%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.756"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("D18").Value = "40.507.07"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("E19").Value = "  -3.35%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -5.03%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.56"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -3.04%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.72"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -4.89%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.02"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("E24").Value = "  -6.26%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -8.48%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.80"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -3.92%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.22%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -4.63%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.22"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -6.20%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.09"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -2.08%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -5.12%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0726"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("E36").Value = "  -2.35%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.04"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -7.13%  "
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("E39").Value = "  -4.95%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.71"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -7.79%  "
$ws.Range("E41").Value = "  -3.45%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.40"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("D43").Value = "1.957.62"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("E44").Value = "  -5.05%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.59"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -6.80%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.33"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -1.34%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -9.73%  "
$ws.Range("D48").Value = "2.599.05"
$ws.Range("E48").Value = "  -4.04%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.92"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -5.09%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.77"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -5.50%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.01"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -4.44%  "
